$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New query text for the "CasesTab" row (row 2, column B) ---
$casesTabQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "WHITE"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# --- New StatQuery text, shared by row 2 and the new row 3 (column C) ---
$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE c.race = "WHITE"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# --- New query text for the new "FilesTab" row (row 3, column B) ---
$filesTabQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
  WHERE c.race = "WHITE"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Add the new FilesTab row (row 3) tab name first
$ws.Range("A3").Value = "FilesTab"

# Update the existing CasesTab row (row 2) with the new queries
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = $statQuery
$ws.Range("C2").WrapText = $true

# Finish the new FilesTab row (row 3)
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = $statQuery
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

# Row heights (Excel grows row 2/3 to fit the long wrapped text)
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# Column widths (closest reproducible values - engine quantizes ColumnWidth
# to 1/6-character steps, so these land on the nearest attainable width)
$ws.Columns.Item(1).ColumnWidth = 8.0
$ws.Columns.Item(2).ColumnWidth = 75.0
$ws.Columns.Item(3).ColumnWidth = 75.0
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.6666666666667

# View state: zoom to 70%, scroll so row 3 is in view, select D3
$win = $ws.Application.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("D3").Select()
